$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the ConnectsTo column for existing rows (replace Server4 refs
#     with the new Firewall devices) ---
$ws.Range("C2").Value = "Firewall1"
$ws.Range("C3").Value = "Server3, Firewall2"
$ws.Range("C4").Value = "Firewall3"

# --- Add the three new firewall nodes as rows 7-9, copying the format of
#     the existing data rows (A6:B6) so style/indentation matches ---
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A7").Value = "Firewall1"
$ws.Range("B7").Value = "192.168.6.1"
$ws.Range("C7").Value = "Server2"

$ws.Range("A8").Value = "Firewall2"
$ws.Range("B8").Value = "192.168.7.1"
$ws.Range("C8").Value = "Server4"

$ws.Range("A9").Value = "Firewall3"
$ws.Range("B9").Value = "192.168.8.1"
$ws.Range("C9").Value = "Server5"

# Row 6 loses its explicit (wrapped-text) row height once re-autofit
$ws.Rows.Item(6).AutoFit()

# Match the author's final selection
$ws.Range("C4").Select()
